# Scheduled runner update: refresh cached marketboard price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-class Leve
# sheets. Values below mirror the latest pull; a handful of rows lose their
# previously-populated NQ/HQ profit cell entirely where the source no longer
# reports that split.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 782.4545000000001
$ws.Range("I8").Value = 782.4545000000001
$ws.Range("K8").Value = 2347.3635
$ws.Range("M8").Value = -2208.3635

# Row 15
$ws.Range("H15").Value = 921.7692
$ws.Range("I15").Value = 921.7692
$ws.Range("K15").Value = 2765.3076
$ws.Range("M15").Value = -2596.3076

# Row 17
$ws.Range("H17").Value = 2973.5
$ws.Range("J17").Value = 2973.5
$ws.Range("L17").Value = 8920.5
$ws.Range("N17").Value = -9256.5

# Row 70
$ws.Range("H70").Value = 2127
$ws.Range("I70").Value = 1599.8
$ws.Range("J70").Value = 3445
$ws.Range("K70").Value = 4799.4
$ws.Range("L70").Value = 10335
$ws.Range("M70").Value = -4529.4
$ws.Range("N70").Value = -10875

# Row 73
$ws.Range("H73").Value = 2127
$ws.Range("I73").Value = 1599.8
$ws.Range("J73").Value = 3445
$ws.Range("K73").Value = 4799.4
$ws.Range("L73").Value = 10335
$ws.Range("M73").Value = -3863.4
$ws.Range("N73").Value = -12207

# Row 76
$ws.Range("H76").Value = 4559.875
$ws.Range("I76").Value = 4332.6665
$ws.Range("K76").Value = 4332.6665
$ws.Range("M76").Value = -4017.6665

# Row 79
$ws.Range("H79").Value = 4559.875
$ws.Range("I79").Value = 4332.6665
$ws.Range("K79").Value = 4332.6665
$ws.Range("M79").Value = -3240.6665

# Row 105
$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988

# Row 112
$ws.Range("H112").Value = 2954.6667
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2954.6667
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 8864.000100000001
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -11080.0001

# Row 121
$ws.Range("H121").Value = 1771.7273
$ws.Range("J121").Value = 1771.7273
$ws.Range("L121").Value = 5315.1819
$ws.Range("N121").Value = -8809.1819

# Row 125
$ws.Range("H125").Value = 2637.4
$ws.Range("I125").Value = 2249.5
$ws.Range("K125").Value = 20245.5
$ws.Range("M125").Value = -17785.5

# Row 135
$ws.Range("H135").Value = 1054.2354
$ws.Range("I135").Value = 970.5
$ws.Range("J135").Value = 2394
$ws.Range("K135").Value = 8734.5
$ws.Range("L135").Value = 21546
$ws.Range("M135").Value = -6199.5
$ws.Range("N135").Value = -26616

# Row 137
$ws.Range("H137").Value = 25643566
$ws.Range("I137").Value = 111112050
$ws.Range("J137").Value = 3022.3
$ws.Range("K137").Value = 333336150
$ws.Range("L137").Value = 9066.900000000001
$ws.Range("M137").Value = -333333600
$ws.Range("N137").Value = -14166.9

# Row 138
$ws.Range("H138").Value = 8736.6875
$ws.Range("I138").Value = 6382.5713
$ws.Range("J138").Value = 10567.667
$ws.Range("K138").Value = 19147.7139
$ws.Range("L138").Value = 31703.001
$ws.Range("M138").Value = -14007.7139
$ws.Range("N138").Value = -41983.001

# Row 139
$ws.Range("H139").Value = 69997
$ws.Range("J139").Value = 69997
$ws.Range("L139").Value = 69997
$ws.Range("N139").Value = -80277

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 11643.3125
$ws.Range("I32").Value = 11643.3125
$ws.Range("K32").Value = 11643.3125
$ws.Range("M32").Value = -11356.3125

# Row 61
$ws.Range("H61").Value = 6464.722
$ws.Range("I61").Value = 5960.375
$ws.Range("K61").Value = 5960.375
$ws.Range("M61").Value = -5748.375

# Row 96
$ws.Range("H96").Value = 16666
$ws.Range("J96").Value = 16666
$ws.Range("L96").Value = 16666
$ws.Range("N96").Value = -22158

# Row 104
$ws.Range("H104").Value = 45414
$ws.Range("J104").Value = 47218.668
$ws.Range("L104").Value = 47218.668
$ws.Range("N104").Value = -54206.668

# Row 105
$ws.Range("H105").Value = 29000
$ws.Range("J105").Value = 29000
$ws.Range("L105").Value = 29000
$ws.Range("N105").Value = -35988

# Row 131
$ws.Range("H131").Value = 80000
$ws.Range("J131").Value = 80000
$ws.Range("L131").Value = 80000
$ws.Range("N131").Value = -90080

# Row 136
$ws.Range("H136").Value = 6464.722
$ws.Range("I136").Value = 5960.375
$ws.Range("K136").Value = 17881.125
$ws.Range("M136").Value = -15331.125

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 350.875
$ws.Range("I5").Value = 95
$ws.Range("J5").Value = 606.75
$ws.Range("K5").Value = 95
$ws.Range("L5").Value = 606.75
$ws.Range("M5").Value = 18
$ws.Range("N5").Value = -832.75

# Row 86
$ws.Range("H86").Value = 30322458
$ws.Range("I86").Value = 29456.857
$ws.Range("J86").Value = 83335210
$ws.Range("K86").Value = 29456.857
$ws.Range("L86").Value = 83335210
$ws.Range("M86").Value = -28333.857
$ws.Range("N86").Value = -83337456

# Row 89
$ws.Range("H89").Value = 30322458
$ws.Range("I89").Value = 29456.857
$ws.Range("J89").Value = 83335210
$ws.Range("K89").Value = 147284.285
$ws.Range("L89").Value = 416676050
$ws.Range("M89").Value = -141668.285
$ws.Range("N89").Value = -416687282

# Row 112
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

# Row 134
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2987.6667
$ws.Range("I31").Value = 2999.9092
$ws.Range("J31").Value = 2954
$ws.Range("K31").Value = 2999.9092
$ws.Range("L31").Value = 2954
$ws.Range("M31").Value = -2704.9092
$ws.Range("N31").Value = -3544

# Row 34
$ws.Range("H34").Value = 2987.6667
$ws.Range("I34").Value = 2999.9092
$ws.Range("J34").Value = 2954
$ws.Range("K34").Value = 2999.9092
$ws.Range("L34").Value = 2954
$ws.Range("M34").Value = -2797.9092
$ws.Range("N34").Value = -3358

# Row 57
$ws.Range("H57").Value = 983.3333
$ws.Range("I57").Value = 983.3333
$ws.Range("K57").Value = 983.3333
$ws.Range("M57").Value = -423.3333

# Row 122
$ws.Range("H122").Value = 5798.625
$ws.Range("I122").Value = 6198.5713
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 18595.7139
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -16145.7139
$ws.Range("N122").Value = -13897

$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Range("H6").Value = 77404.08
$ws.Range("I6").Value = 83521.086
$ws.Range("K6").Value = 250563.258
$ws.Range("M6").Value = -250450.258

# Row 20
$ws.Range("H20").Value = 3500
$ws.Range("I20").Value = 3000
$ws.Range("J20").Value = 4000
$ws.Range("K20").Value = 9000
$ws.Range("L20").Value = 12000
$ws.Range("M20").Value = -8773
$ws.Range("N20").Value = -12454

# Row 33
$ws.Range("H33").Value = 282.875
$ws.Range("I33").Value = 277.57144
$ws.Range("K33").Value = 1665.42864
$ws.Range("M33").Value = -1382.42864

# Row 99
$ws.Range("H99").Value = 2902
$ws.Range("I99").Value = 2594.8333
$ws.Range("K99").Value = 7784.499899999999
$ws.Range("M99").Value = -5538.499899999999

# Row 128
$ws.Range("H128").Value = 199999
$ws.Range("I128").Value = 199999
$ws.Range("K128").Value = 599997
$ws.Range("M128").Value = -595017

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 254.71428
$ws.Range("I2").Value = 163.55556
$ws.Range("J2").Value = 418.8
$ws.Range("K2").Value = 163.55556
$ws.Range("L2").Value = 418.8
$ws.Range("M2").Value = -50.55556000000001
$ws.Range("N2").Value = -644.8

# Row 113
$ws.Range("H113").Value = 1840
$ws.Range("I113").Value = 1820
$ws.Range("K113").Value = 1820
$ws.Range("M113").Value = 350

# Row 122
$ws.Range("H122").Value = 5149.2144
$ws.Range("I122").Value = 4967.409
$ws.Range("J122").Value = 5815.8335
$ws.Range("K122").Value = 14902.227
$ws.Range("L122").Value = 17447.5005
$ws.Range("M122").Value = -12452.227
$ws.Range("N122").Value = -22347.5005

$ws = $wb.Worksheets.Item("LTW")
# Row 43
$ws.Range("H43").Value = 999
$ws.Range("I43").Value = 999
$ws.Range("K43").Value = 999
$ws.Range("M43").Value = -806

# Row 45
$ws.Range("H45").Value = 9999.5
$ws.Range("I45").Value = 9999.5
$ws.Range("K45").Value = 9999.5
$ws.Range("M45").Value = -9592.5

# Row 94
$ws.Range("H94").Value = 76999.75
$ws.Range("J94").Value = 76999.75
$ws.Range("L94").Value = 76999.75
$ws.Range("N94").Value = -78351.75

# Row 105
$ws.Range("H105").Value = 20000
$ws.Range("J105").Value = 20000
$ws.Range("L105").Value = 20000
$ws.Range("N105").Value = -26988

# Row 106
$ws.Range("H106").Value = 16680
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 16680
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 16680
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -19204

$ws = $wb.Worksheets.Item("WVR")
# Row 101
$ws.Range("H101").Value = 28404.25
$ws.Range("J101").Value = 28404.25
$ws.Range("L101").Value = 28404.25
$ws.Range("N101").Value = -34894.25

# Row 104
$ws.Range("H104").Value = 10000
$ws.Range("J104").Value = 10000
$ws.Range("L104").Value = 10000
$ws.Range("N104").Value = -16988

# Row 132
$ws.Range("H132").Value = 500001180
$ws.Range("I132").Value = 2345
$ws.Range("K132").Value = 7035
$ws.Range("M132").Value = -4505
